# Update countries & provincias Spain
# Refreshes the COVID stats table: updates the "last updated" timestamp,
# re-applies the latest per-country figures (cols B-H), and fixes the
# country label for rows whose rank order changed now that the figures
# moved (the source list is kept sorted by "Casos totales" descending).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." banner in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 01:54"

# Row => country label (re-asserted even where unchanged, harmless) + new B..H figures
$changes = @(
    @{ Row=4;   Name="Estados Unidos";               B=6297654; C=40083; D=3537272; E=2570453; F=0; G=1029; H=189929 }
    @{ Row=5;   Name="Brasil";                        B=4001422; C=48632; D=3210405; E=667118;  F=0; G=1218; H=123899 }
    @{ Row=8;   Name="Peru";                          B=663437;  C=6308;  D=480177;  E=154001;  F=0; G=191;  H=29259 }
    @{ Row=13;  Name="Argentina";                     B=439172;  C=10933; D=315530;  E=114524;  F=0; G=199;  H=9118 }
    @{ Row=36;  Name="Panama";                        B=94084;   C=532;   D=67487;   E=24567;   F=0; G=12;   H=2030 }
    @{ Row=52;  Name="Singapur";                      B=56860;   C=8;     D=55891;   E=942;     F=0; G=0;    H=27 }
    @{ Row=57;  Name="Argelia";                       B=45158;   C=325;   D=31746;   E=11889;   F=0; G=5;    H=1523 }
    @{ Row=58;  Name="Ghana";                         B=44658;   C=198;   D=43478;   E=904;     F=0; G=0;    H=276 }
    @{ Row=74;  Name="Chequia";                       B=25773;   C=656;   D=18326;   E=7022;    F=0; G=0;    H=425 }
    @{ Row=99;  Name="Gabon";                         B=8538;    C=5;     D=7335;    E=1150;    F=0; G=0;    H=53 }
    @{ Row=105; Name="Luxemburgo";                    B=6745;    C=43;    D=0;       E=0;       F=0; G=0;    H=124 }
    @{ Row=106; Name="Zimbabue";                      B=6638;    C=79;    D=5250;    E=1182;    F=0; G=3;    H=206 }
    @{ Row=107; Name="Hungria";                       B=6622;    C=365;   D=3903;    E=2100;    F=0; G=3;    H=619 }
    @{ Row=108; Name="Malaui";                        B=5579;    C=3;     D=3500;    E=1904;    F=0; G=0;    H=175 }
    @{ Row=110; Name="Montenegro";                    B=5019;    C=102;   D=4104;    E=811;     F=0; G=2;    H=104 }
    @{ Row=111; Name="Guinea Ecuatorial";              B=4965;    C=0;     D=3884;    E=998;     F=0; G=0;    H=83 }
    @{ Row=113; Name="Republica de Africa Central";    B=4712;    C=1;     D=1803;    E=2847;    F=0; G=0;    H=62 }
    @{ Row=118; Name="Surinam";                        B=4149;    C=60;    D=3272;    E=805;     F=0; G=0;    H=72 }
    @{ Row=119; Name="Cuba";                           B=4126;    C=61;    D=3458;    E=570;     F=0; G=3;    H=98 }
    @{ Row=120; Name="Mozambique";                     B=4117;    C=78;    D=2170;    E=1922;    F=0; G=2;    H=25 }
    @{ Row=141; Name="Benin";                          B=2194;    C=49;    D=1738;    E=416;     F=0; G=0;    H=40 }
    @{ Row=142; Name="Jordania";                       B=2161;    C=64;    D=1610;    E=536;     F=0; G=0;    H=15 }
    @{ Row=151; Name="Uruguay";                        B=1626;    C=15;    D=1433;    E=149;     F=0; G=0;    H=44 }
    @{ Row=167; Name="Santo Tome y Principe";           B=896;     C=0;     D=855;     E=26;      F=0; G=0;    H=15 }
    @{ Row=171; Name="Polinesia Francesa";              B=622;     C=49;    D=336;     E=286;     F=0; G=0;    H=0 }
    @{ Row=176; Name="Papua Nueva Guinea";               B=471;     C=11;    D=232;     E=234;     F=0; G=0;    H=5 }
    @{ Row=177; Name="Burundi";                          B=448;     C=3;     D=357;     E=90;      F=0; G=0;    H=1 }
)

foreach ($chg in $changes) {
    $r = $chg.Row
    $ws.Cells.Item($r, 1).Value = $chg.Name
    $ws.Cells.Item($r, 2).Value = $chg.B
    $ws.Cells.Item($r, 3).Value = $chg.C
    $ws.Cells.Item($r, 4).Value = $chg.D
    $ws.Cells.Item($r, 5).Value = $chg.E
    $ws.Cells.Item($r, 6).Value = $chg.F
    $ws.Cells.Item($r, 7).Value = $chg.G
    $ws.Cells.Item($r, 8).Value = $chg.H
}
